$d = $word.ActiveDocument

# Map of exact old paragraph text -> new Swahili (Kenya) translation.
# Each key occurs exactly once as the full text of a paragraph in the
# document (verified against the source before writing this script), so
# matching on the complete paragraph text (rather than a blind Find/Replace)
# avoids accidentally touching similarly-worded longer strings such as
# "Solution of the first experiment" or "Assist the process, provoke
# thoughts (explanations from 05:26 - 06:06)".
$translations = [ordered]@{
    "Video Title"                           = "Kichwa cha Video"
    "Topic"                                 = "Mada"
    "Aim(s)"                                = "Malengo"
    "Length"                                = "Urefu"
    "Camp Location"                         = "Mahali pa Kambi"
    "Facilitators"                          = "Wawezeshaji"
    "N. of students"                        = "N. ya wanafunzi"
    "Date"                                  = "Tarehe"
    "Resources"                             = "Rasilimali"
    "needed"                                = "inahitajika"
    "Preparations"                          = "Maandalizi"
    "Video time"                            = "Muda wa video"
    "What facilitator does"                 = "Mwezeshaji anafanya nini"
    "What learners do"                      = "Wanachofanya wanafunzi"
    "General VMC Video Introduction"        = "Utangulizi Mkuu wa Video ya VMC"
    "Video Introduction"                    = "Utangulizi wa Video"
    "Assist the process, provoke thoughts"  = "Kusaidia mchakato, kuchochea mawazo"
    "Solution"                              = "Suluhisho"
}

foreach ($p in $d.Paragraphs) {
    $full = $p.Range.Text
    # Paragraphs end with a paragraph mark (Cr, 13); ones inside table
    # cells are further followed by the cell mark (Bell, 7). Strip both
    # so the comparison is against the paragraph's visible text only.
    $plain = $full.TrimEnd([char]13, [char]7)

    if ($translations.Contains($plain)) {
        $newText = $translations[$plain]
        $r = $p.Range
        # Shrink the range so it covers only the visible text, leaving the
        # paragraph mark (and cell mark, if any) untouched; this preserves
        # paragraph/run formatting exactly as Word would when you retype
        # the contents of a run in place.
        $r.MoveEnd(1, -1) | Out-Null
        $r.Text = $newText
    }
}

# Document default language: Swahili (Tanzania) -> Swahili (Kenya).
# This mirrors the one non-text change in the diff (the rPrDefault language
# tag). The Word object model exposes the document's default/base language
# through the "Normal" style, so update it there.
$normal = $d.Styles(-1)
$normal.LanguageID = "sw-KE"
